# Auto update Excel log
# Appends new mmWave sensor log rows (184-204) to the "mmWave" worksheet,
# extending the used range from A1:F183 to A1:F204.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

# Each entry: Row, Date, Timestamp, Hour, Location, Value, Status
$data = @(
    @("184", "2026-01-28", "17:50:42", "17:00", "Living Room", "PRESENCE",    "Active"),
    @("185", "2026-01-28", "17:50:43", "17:00", "Living Room", "PRESENCE",    "Active"),
    @("186", "2026-01-28", "17:50:45", "17:00", "Living Room", "PRESENCE",    "Active"),
    @("187", "2026-01-28", "17:50:49", "17:00", "Living Room", "PRESENCE",    "Active"),
    @("188", "2026-01-28", "17:50:52", "17:00", "Living Room", "PRESENCE",    "Active"),
    @("189", "2026-01-28", "17:50:55", "17:00", "Living Room", "PRESENCE",    "Active"),
    @("190", "2026-01-28", "17:50:57", "17:00", "Living Room", "PRESENCE",    "Active"),
    @("191", "2026-01-28", "17:51:01", "17:00", "Living Room", "PRESENCE",    "Active"),
    @("192", "2026-01-28", "17:51:04", "17:00", "Living Room", "PRESENCE",    "Active"),
    @("193", "2026-01-28", "17:51:06", "17:00", "Living Room", "PRESENCE",    "Active"),
    @("194", "2026-01-28", "17:51:10", "17:00", "Living Room", "PRESENCE",    "Active"),
    @("195", "2026-01-28", "17:51:12", "17:00", "Living Room", "PRESENCE",    "Active"),
    @("196", "2026-01-28", "17:51:15", "17:00", "Living Room", "PRESENCE",    "Active"),
    @("197", "2026-01-28", "17:51:19", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("198", "2026-01-28", "17:51:22", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("199", "2026-01-28", "17:51:24", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("200", "2026-01-28", "17:51:28", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("201", "2026-01-28", "17:51:31", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("202", "2026-01-28", "17:51:34", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("203", "2026-01-28", "17:51:36", "17:00", "Living Room", "NO_PRESENCE", "Inactive"),
    @("204", "2026-01-28", "17:51:40", "17:00", "Living Room", "NO_PRESENCE", "Inactive")
)

foreach ($row in $data) {
    $r = $row[0]

    # Column A holds an ISO-like date string ("2026-01-28"). Excel would
    # otherwise auto-convert this into a real date serial number, so a
    # leading apostrophe forces it to be kept as literal text, matching
    # the rest of the column. Re-applying the "Normal" style afterwards
    # clears the quote-prefix formatting flag that the apostrophe leaves
    # behind, so the cell ends up with the default (unstyled) format.
    $ws.Range("A$r").Value = "'" + $row[1]
    $ws.Range("A$r").Style = "Normal"

    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
    $ws.Range("F$r").Value = $row[6]
}

Write-Host "Appended $($data.Count) rows to mmWave sheet (184-204)."
